$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 (Marking): Right count decreases, Wrong count becomes more negative
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -2

# Row 12 (Total): total marks corrected, reflecting the corrected marking row
$ws.Range("B12").Value = 104
$ws.Range("C12").Value = -4
$ws.Range("E12").Value = "100 / 112"
